$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# for rows 2-13, per the weekly data re-shuffle.

$rows = @{
    2  = @{ D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
    3  = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    4  = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    5  = @{ D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    6  = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    7  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 }
    8  = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 }
    9  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
    10 = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    11 = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    12 = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    13 = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
}
